$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# New "Link" entry for SP (row 6) and new "What have we done?" status column (L)
$ws.Range("E6").Value = "http://www.sp.se/en/workatsp/students/masterthesis/Sidor/STRATEGIESANDMETHODSFORSUSTAINABLEINNOVATIONANDTRANSFORMATION.aspx"
$ws.Range("L1").Value = "What have we done?"
$ws.Range("L6").Value = "Sent mail, waiting for answer."
$ws.Range("L11").Value = "Sent mail, waiting for answer."

# Column K (Deadline) widens to fit its longest entry once the sheet is touched again
$ws.Columns("K:K").AutoFit() | Out-Null

# Leave the selection where the author left it when done editing
$ws.Range("L8").Select() | Out-Null
